# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    existing "2022-Q1" sheet) and populate it with the fund-holding detail
#    rows for that quarter.
# 2. Insert a new summary row at the top of the "总计" sheet for "2022-Q4"
#    (15 holdings, 0.76 亿元 market value), pushing the existing "2022-Q1"
#    and "2021-Q2" summary rows down by one row each.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: stamp a range with the same look as the workbook's existing
# "header / index column" style (bold font, centered, thin box border) —
# this is style index 2 in the original file, applied to A-column index
# cells and header cells on every detail sheet.
# ---------------------------------------------------------------------
function Set-HeaderStyle($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108
    $rng.VerticalAlignment = -4160
    $rng.Borders.LineStyle = 1
}

# =======================================================================
# Step 1: create the new "2022-Q4" sheet, placed before "2022-Q1"
# =======================================================================
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q4Sheet = $wb.Worksheets.Add($q1Sheet, $null)
$q4Sheet.Name = "2022-Q4"

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"
Set-HeaderStyle $q4Sheet.Range("B1:H1")

# Fund holding detail rows (row 2 .. row 16)
$data = @(
    @("010874", "泰康品质生活混合A", "6.77", "84.20", "2.99", "0.2024", 9),
    @("160211", "国泰中小盘成长混合（LOF）", "6.59", "90.98", "1.97", "0.1298", 9),
    @("010875", "泰康品质生活混合C", "3.33", "84.20", "2.99", "0.0996", 9),
    @("005416", "鹏华尊惠18个月定期开放混合A", "2.83", "39.19", "1.87", "0.0529", 6),
    @("009668", "鹏华安庆混合C", "2.79", "39.85", "1.68", "0.0469", 5),
    @("003166", "鹏华弘嘉灵活配置混合C", "0.92", "91.65", "4.40", "0.0405", 4),
    @("009667", "鹏华安庆混合A", "2.34", "39.85", "1.68", "0.0393", 5),
    @("003165", "鹏华弘嘉灵活配置混合A", "0.72", "91.65", "4.40", "0.0317", 4),
    @("011573", "鹏华安荣混合C", "1.87", "39.61", "1.54", "0.0288", 6),
    @("011572", "鹏华安荣混合A", "1.54", "39.61", "1.54", "0.0237", 6),
    @("009231", "鹏华安和混合C", "1.54", "38.20", "1.53", "0.0236", 6),
    @("009230", "鹏华安和混合A", "1.37", "38.20", "1.53", "0.0210", 6),
    @("005417", "鹏华尊惠18个月定期开放混合C", "0.36", "39.19", "1.87", "0.0067", 6),
    @("008324", "宝盈祥利稳健配置混合A", "0.53", "31.08", "1.19", "0.0063", 10),
    @("008325", "宝盈祥利稳健配置混合C", "0.29", "31.08", "1.19", "0.0035", 10)
)

$row = 2
foreach ($item in $data) {
    $idx = $row - 2

    $aCell = $q4Sheet.Range("A$row")
    $aCell.Value = $idx
    Set-HeaderStyle $aCell

    # Text-typed columns (must keep leading zeros / exact decimal text)
    $q4Sheet.Range("B$row").NumberFormat = "@"
    $q4Sheet.Range("B$row").Value = $item[0]
    $q4Sheet.Range("C$row").NumberFormat = "@"
    $q4Sheet.Range("C$row").Value = $item[1]
    $q4Sheet.Range("D$row").NumberFormat = "@"
    $q4Sheet.Range("D$row").Value = $item[2]
    $q4Sheet.Range("E$row").NumberFormat = "@"
    $q4Sheet.Range("E$row").Value = $item[3]
    $q4Sheet.Range("F$row").NumberFormat = "@"
    $q4Sheet.Range("F$row").Value = $item[4]
    $q4Sheet.Range("G$row").NumberFormat = "@"
    $q4Sheet.Range("G$row").Value = $item[5]

    # Numeric column
    $q4Sheet.Range("H$row").Value = $item[6]

    $row = $row + 1
}

# =======================================================================
# Step 2: update the "总计" (summary) sheet — insert the new 2022-Q4 row
# at the top of the data and push everything else down.
# =======================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$a2 = $totalSheet.Range("A2")
$a2.Value = 0
Set-HeaderStyle $a2
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 15
$totalSheet.Range("D2").Value = 0.76

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.05

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q2"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.06
